# Daily attendance processing - 2026-01-27 22:00:31
# Normalize the "Recorded By" (column G) value ordering for affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 1) {
    $lastRow = 157
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Text

    if ($current -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($current -eq "backup@backdoor.com, System, system") {
        $cell.Value2 = "system, backup@backdoor.com, System"
    }
    elseif ($current -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, admin@admin.com"
    }
}
